$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update task rows 2..34 (A: id, B: task name, C: resources, E: start, F: end, G: outline level) ---
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "Davide"
$ws.Cells.Item(2, 3).Value = $null
$ws.Cells.Item(2, 5).Value = 42660.333333333336
$ws.Cells.Item(2, 6).Value = 42723.666666666664
$ws.Cells.Item(2, 7).Value = 1

$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "Meetings with client"
$ws.Cells.Item(3, 3).Value = "Mario, Moreno"
$ws.Cells.Item(3, 5).Value = 42660.333333333336
$ws.Cells.Item(3, 6).Value = 42662.666666666664
$ws.Cells.Item(3, 7).Value = 2

$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "Identification of goals and stakeholders"
$ws.Cells.Item(4, 3).Value = "Mario, Moreno"
$ws.Cells.Item(4, 5).Value = 42660.333333333336
$ws.Cells.Item(4, 6).Value = 42660.666666666664
$ws.Cells.Item(4, 7).Value = 3

$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "Requirements elicitation"
$ws.Cells.Item(5, 3).Value = "Mario, Moreno"
$ws.Cells.Item(5, 5).Value = 42661.333333333336
$ws.Cells.Item(5, 6).Value = 42662.666666666664
$ws.Cells.Item(5, 7).Value = 3

$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "Modelization of the World and the Machine"
$ws.Cells.Item(6, 3).Value = "Mario, Moreno"
$ws.Cells.Item(6, 5).Value = 42663.333333333336
$ws.Cells.Item(6, 6).Value = 42668.666666666664
$ws.Cells.Item(6, 7).Value = 2

$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "Identification of domain assumptions"
$ws.Cells.Item(7, 3).Value = "Mario, Moreno"
$ws.Cells.Item(7, 5).Value = 42663.333333333336
$ws.Cells.Item(7, 6).Value = 42668.666666666664
$ws.Cells.Item(7, 7).Value = 3

$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "Identification of system goals"
$ws.Cells.Item(8, 3).Value = "Mario, Moreno"
$ws.Cells.Item(8, 5).Value = 42663.333333333336
$ws.Cells.Item(8, 6).Value = 42668.666666666664
$ws.Cells.Item(8, 7).Value = 3

$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "Identification of requirements"
$ws.Cells.Item(9, 3).Value = "Mario, Moreno"
$ws.Cells.Item(9, 5).Value = 42663.333333333336
$ws.Cells.Item(9, 6).Value = 42668.666666666664
$ws.Cells.Item(9, 7).Value = 3

$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "Writing scenarios"
$ws.Cells.Item(10, 3).Value = "Moreno"
$ws.Cells.Item(10, 5).Value = 42669.333333333336
$ws.Cells.Item(10, 6).Value = 42669.6666666088
$ws.Cells.Item(10, 7).Value = 2

$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "Identification of needed data"
$ws.Cells.Item(11, 3).Value = $null
$ws.Cells.Item(11, 5).Value = 42670.333333333336
$ws.Cells.Item(11, 6).Value = 42671.666666666664
$ws.Cells.Item(11, 7).Value = 2

$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "Identification of use cases"
$ws.Cells.Item(12, 3).Value = "Moreno"
$ws.Cells.Item(12, 5).Value = 42671.333333333336
$ws.Cells.Item(12, 6).Value = 42671.666666666664
$ws.Cells.Item(12, 7).Value = 2

$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = "In progress meeting with client"
$ws.Cells.Item(13, 3).Value = "Mario, Moreno"
$ws.Cells.Item(13, 5).Value = 42676.333333333336
$ws.Cells.Item(13, 6).Value = 42676.666666666664
$ws.Cells.Item(13, 7).Value = 2

$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = "Document refinement"
$ws.Cells.Item(14, 3).Value = "Mario, Moreno"
$ws.Cells.Item(14, 5).Value = 42677.333333333336
$ws.Cells.Item(14, 6).Value = 42687.666666666664
$ws.Cells.Item(14, 7).Value = 2

$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = "Requirements refinement"
$ws.Cells.Item(15, 3).Value = "Mario, Moreno"
$ws.Cells.Item(15, 5).Value = 42677.333333333336
$ws.Cells.Item(15, 6).Value = 42681.666666666664
$ws.Cells.Item(15, 7).Value = 3

$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = "Data model refinement"
$ws.Cells.Item(16, 3).Value = $null
$ws.Cells.Item(16, 5).Value = 42677.333333333336
$ws.Cells.Item(16, 6).Value = 42678.666666666664
$ws.Cells.Item(16, 7).Value = 3

$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = "Alloy modelization"
$ws.Cells.Item(17, 3).Value = "Mario"
$ws.Cells.Item(17, 5).Value = 42677.333333333336
$ws.Cells.Item(17, 6).Value = 42681.666666666664
$ws.Cells.Item(17, 7).Value = 3

$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = "Document revision"
$ws.Cells.Item(18, 3).Value = "Mario, Moreno"
$ws.Cells.Item(18, 5).Value = 42682.333333333336
$ws.Cells.Item(18, 6).Value = 42685.666666666664
$ws.Cells.Item(18, 7).Value = 2

$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = "Architecture draft"
$ws.Cells.Item(19, 3).Value = "Mario, Moreno"
$ws.Cells.Item(19, 5).Value = 42688.333333333336
$ws.Cells.Item(19, 6).Value = 42690.666666666664
$ws.Cells.Item(19, 7).Value = 2

$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = "High level system view"
$ws.Cells.Item(20, 3).Value = "Mario, Moreno"
$ws.Cells.Item(20, 5).Value = 42688.333333333336
$ws.Cells.Item(20, 6).Value = 42690.666666666664
$ws.Cells.Item(20, 7).Value = 3

$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = "Definition of system boundaries"
$ws.Cells.Item(21, 3).Value = "Mario, Moreno"
$ws.Cells.Item(21, 5).Value = 42688.333333333336
$ws.Cells.Item(21, 6).Value = 42690.666666666664
$ws.Cells.Item(21, 7).Value = 3

$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = "Main architecture decisions"
$ws.Cells.Item(22, 3).Value = "Mario, Moreno"
$ws.Cells.Item(22, 5).Value = 42688.333333333336
$ws.Cells.Item(22, 6).Value = 42690.666666666664
$ws.Cells.Item(22, 7).Value = 3

$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = "Meeting with clients"
$ws.Cells.Item(23, 3).Value = "Mario, Moreno"
$ws.Cells.Item(23, 5).Value = 42691.333333333336
$ws.Cells.Item(23, 6).Value = 42691.666666666664
$ws.Cells.Item(23, 7).Value = 2

$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(24, 2).Value = "Refining architecture choices"
$ws.Cells.Item(24, 3).Value = "Mario, Moreno"
$ws.Cells.Item(24, 5).Value = 42692.333333333336
$ws.Cells.Item(24, 6).Value = 42695.666666666664
$ws.Cells.Item(24, 7).Value = 2

$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 2).Value = "Main components component diagram"
$ws.Cells.Item(25, 3).Value = "Mario, Moreno"
$ws.Cells.Item(25, 5).Value = 42696.333333333336
$ws.Cells.Item(25, 6).Value = 42697.666666666664
$ws.Cells.Item(25, 7).Value = 2

$ws.Cells.Item(26, 1).Value = 25
$ws.Cells.Item(26, 2).Value = "ER DB Model"
$ws.Cells.Item(26, 3).Value = $null
$ws.Cells.Item(26, 5).Value = 42698.333333333336
$ws.Cells.Item(26, 6).Value = 42698.666666666664
$ws.Cells.Item(26, 7).Value = 2

$ws.Cells.Item(27, 1).Value = 26
$ws.Cells.Item(27, 2).Value = "Sequence diagrams"
$ws.Cells.Item(27, 3).Value = "Moreno"
$ws.Cells.Item(27, 5).Value = 42699.333333333336
$ws.Cells.Item(27, 6).Value = 42704.666666666664
$ws.Cells.Item(27, 7).Value = 2

$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = "Algorithms desing"
$ws.Cells.Item(28, 3).Value = $null
$ws.Cells.Item(28, 5).Value = 42705.333333333336
$ws.Cells.Item(28, 6).Value = 42705.666666666664
$ws.Cells.Item(28, 7).Value = 2

$ws.Cells.Item(29, 1).Value = 28
$ws.Cells.Item(29, 2).Value = "Revision and requirement traceability"
$ws.Cells.Item(29, 3).Value = "Mario, Moreno"
$ws.Cells.Item(29, 5).Value = 42706.333333333336
$ws.Cells.Item(29, 6).Value = 42706.666666666664
$ws.Cells.Item(29, 7).Value = 2

$ws.Cells.Item(30, 1).Value = 29
$ws.Cells.Item(30, 2).Value = "Integration test strategy"
$ws.Cells.Item(30, 3).Value = "Mario, Moreno"
$ws.Cells.Item(30, 5).Value = 42716.333333333336
$ws.Cells.Item(30, 6).Value = 42716.666666666664
$ws.Cells.Item(30, 7).Value = 2

$ws.Cells.Item(31, 1).Value = 30
$ws.Cells.Item(31, 2).Value = "Definition of precedences"
$ws.Cells.Item(31, 3).Value = "Mario, Moreno"
$ws.Cells.Item(31, 5).Value = 42716.333333333336
$ws.Cells.Item(31, 6).Value = 42716.666666666664
$ws.Cells.Item(31, 7).Value = 2

$ws.Cells.Item(32, 1).Value = 31
$ws.Cells.Item(32, 2).Value = "Integration test description"
$ws.Cells.Item(32, 3).Value = $null
$ws.Cells.Item(32, 5).Value = 42717.333333333336
$ws.Cells.Item(32, 6).Value = 42719.666666666664
$ws.Cells.Item(32, 7).Value = 2

$ws.Cells.Item(33, 1).Value = 32
$ws.Cells.Item(33, 2).Value = "Data required"
$ws.Cells.Item(33, 3).Value = $null
$ws.Cells.Item(33, 5).Value = 42720.333333333336
$ws.Cells.Item(33, 6).Value = 42720.666666666664
$ws.Cells.Item(33, 7).Value = 2

$ws.Cells.Item(34, 1).Value = 33
$ws.Cells.Item(34, 2).Value = "Document revision"
$ws.Cells.Item(34, 3).Value = "Mario, Moreno"
$ws.Cells.Item(34, 5).Value = 42720.333333333336
$ws.Cells.Item(34, 6).Value = 42723.666666666664
$ws.Cells.Item(34, 7).Value = 2

# --- Duration formula (col D), extended across the new range ---
$ws.Range("D2:D34").Formula = '=CONCATENATE(NETWORKDAYS(E2,F2),"g")'

# --- New trailing K/L placeholder cells (date-formatted, empty) for rows 30-34 ---
$ws.Range("E2").Copy() | Out-Null
foreach ($r in 30..34) {
  $ws.Range("K" + $r).PasteSpecial(-4122) | Out-Null
  $ws.Range("L" + $r).PasteSpecial(-4122) | Out-Null
}
$ws.Range("K30:L34").ClearContents()
$excel.CutCopyMode = $false

# --- Selection moves to F36 ---
$ws.Range("F36").Select() | Out-Null
